$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$top = $s.Shapes.Item(1)
$grp = $top.GroupItems

$sh = $grp.Item(3)
$sh.Left = 158.89110236220472
$sh.Top = 125.24543307086614
$sh.Width = 544.8168503937007
$sh.Height = 319.6067716535433
$sh = $grp.Item(4)
$sh.Left = 234.05779527559056
$sh.Top = 125.24543307086614
$sh.Width = 0.0
$sh.Height = 319.6067716535433
$sh = $grp.Item(5)
$sh.Left = 384.3911023622047
$sh.Top = 125.24543307086614
$sh.Width = 0.0
$sh.Height = 319.6067716535433
$sh = $grp.Item(6)
$sh.Left = 534.7244094488188
$sh.Top = 125.24543307086614
$sh.Width = 0.0
$sh.Height = 319.6067716535433
$sh = $grp.Item(7)
$sh.Left = 685.0577165354331
$sh.Top = 125.24543307086614
$sh.Width = 0.0
$sh.Height = 319.6067716535433
$sh = $grp.Item(8)
$sh.Left = 158.89110236220472
$sh.Top = 413.9225196850394
$sh.Width = 544.8168503937007
$sh.Height = 0.0
$sh = $grp.Item(9)
$sh.Left = 158.89110236220472
$sh.Top = 362.37307086614175
$sh.Width = 544.8168503937007
$sh.Height = 0.0
$sh = $grp.Item(10)
$sh.Left = 158.89110236220472
$sh.Top = 310.8235433070866
$sh.Width = 544.8168503937007
$sh.Height = 0.0
$sh = $grp.Item(11)
$sh.Left = 158.89110236220472
$sh.Top = 259.27409448818895
$sh.Width = 544.8168503937007
$sh.Height = 0.0
$sh = $grp.Item(12)
$sh.Left = 158.89110236220472
$sh.Top = 207.72456692913386
$sh.Width = 544.8168503937007
$sh.Height = 0.0
$sh = $grp.Item(13)
$sh.Left = 158.89110236220472
$sh.Top = 156.1751181102362
$sh.Width = 544.8168503937007
$sh.Height = 0.0
$sh = $grp.Item(14)
$sh.Left = 158.89110236220472
$sh.Top = 125.24543307086614
$sh.Width = 0.0
$sh.Height = 319.6067716535433
$sh = $grp.Item(15)
$sh.Left = 309.2244094488189
$sh.Top = 125.24543307086614
$sh.Width = 0.0
$sh.Height = 319.6067716535433
$sh = $grp.Item(16)
$sh.Left = 459.5577165354331
$sh.Top = 125.24543307086614
$sh.Width = 0.0
$sh.Height = 319.6067716535433
$sh = $grp.Item(17)
$sh.Left = 609.8910236220472
$sh.Top = 125.24543307086614
$sh.Width = 0.0
$sh.Height = 319.6067716535433
$sh = $grp.Item(18)
$sh.Left = 158.89110236220472
$sh.Top = 132.97787401574803
$sh.Width = 534.1341732283464
$sh.Height = 46.394488188976375
$sh = $grp.Item(19)
$sh.Left = 158.89110236220472
$sh.Top = 184.52732283464567
$sh.Width = 461.5232283464567
$sh.Height = 46.394488188976375
$sh = $grp.Item(20)
$sh.Left = 158.89110236220472
$sh.Top = 236.07677165354332
$sh.Width = 331.7855905511811
$sh.Height = 46.394488188976375
$sh = $grp.Item(21)
$sh.Left = 158.89110236220472
$sh.Top = 287.62629921259844
$sh.Width = 315.3992125984252
$sh.Height = 46.394488188976375
$sh = $grp.Item(22)
$sh.Left = 158.89110236220472
$sh.Top = 339.1757480314961
$sh.Width = 124.32559055118111
$sh.Height = 46.394488188976375
$sh = $grp.Item(23)
$sh.Left = 158.89110236220472
$sh.Top = 390.7252755905512
$sh.Width = 9.470944881889764
$sh.Height = 46.394488188976375
$sh = $grp.Item(24)
$sh.Left = 634.5455118110236
$sh.Top = 151.00149606299212
$sh = $grp.Item(25)
$sh.Left = 636.9128346456692
$sh.Top = 154.2024409448819
$sh = $grp.Item(26)
$sh.Left = 636.9128346456692
$sh.Top = 155.6028346456693
$sh = $grp.Item(27)
$sh.Left = 643.0146456692913
$sh.Top = 155.53614173228345
$sh = $grp.Item(28)
$sh.Left = 644.4817322834646
$sh.Top = 156.1363779527559
$sh = $grp.Item(29)
$sh.Left = 646.282283464567
$sh.Top = 156.1363779527559
$sh = $grp.Item(30)
$sh.Left = 648.0494488188976
$sh.Top = 156.1363779527559
$sh = $grp.Item(31)
$sh.Left = 644.248346456693
$sh.Top = 150.9348031496063
$sh = $grp.Item(32)
$sh.Left = 643.4814173228347
$sh.Top = 151.00149606299212
$sh = $grp.Item(33)
$sh.Left = 643.0146456692913
$sh.Top = 153.46889763779527
$sh = $grp.Item(34)
$sh.Left = 647.7493700787402
$sh.Top = 153.46889763779527
$sh = $grp.Item(35)
$sh.Left = 651.2003937007875
$sh.Top = 151.84338582677165
$sh = $grp.Item(36)
$sh.Left = 654.8014173228346
$sh.Top = 151.00149606299212
$sh = $grp.Item(37)
$sh.Left = 657.1688188976378
$sh.Top = 154.2024409448819
$sh = $grp.Item(38)
$sh.Left = 657.1688188976378
$sh.Top = 155.6028346456693
$sh = $grp.Item(39)
$sh.Left = 663.4706299212598
$sh.Top = 150.96818897637795
$sh = $grp.Item(40)
$sh.Left = 663.4373228346457
$sh.Top = 152.601968503937
$sh = $grp.Item(41)
$sh.Left = 664.9377165354331
$sh.Top = 156.8032283464567
$sh = $grp.Item(42)
$sh.Left = 671.8730708661417
$sh.Top = 151.03488188976377
$sh = $grp.Item(43)
$sh.Left = 672.3732283464567
$sh.Top = 153.00212598425196
$sh = $grp.Item(44)
$sh.Left = 673.0400787401575
$sh.Top = 153.60228346456694
$sh = $grp.Item(45)
$sh.Left = 673.3735433070866
$sh.Top = 153.86897637795275
$sh = $grp.Item(46)
$sh.Left = 676.3077165354331
$sh.Top = 156.00299212598426
$sh = $grp.Item(47)
$sh.Left = 565.0521259842519
$sh.Top = 202.5843307086614
$sh = $grp.Item(48)
$sh.Left = 564.9854330708662
$sh.Top = 204.81834645669292
$sh = $grp.Item(49)
$sh.Left = 569.6867716535434
$sh.Top = 202.61771653543306
$sh = $grp.Item(50)
$sh.Left = 570.753779527559
$sh.Top = 204.61826771653543
$sh = $grp.Item(51)
$sh.Left = 567.1526771653544
$sh.Top = 202.61771653543306
$sh = $grp.Item(52)
$sh.Left = 566.9859842519685
$sh.Top = 205.78527559055118
$sh = $grp.Item(53)
$sh.Left = 566.7192913385827
$sh.Top = 206.85228346456694
$sh = $grp.Item(54)
$sh.Left = 573.6879527559055
$sh.Top = 202.5843307086614
$sh = $grp.Item(55)
$sh.Left = 585.458031496063
$sh.Top = 207.21905511811025
$sh = $grp.Item(56)
$sh.Left = 588.1922047244094
$sh.Top = 207.4524409448819
$sh = $grp.Item(57)
$sh.Left = 584.0910236220473
$sh.Top = 207.48574803149606
$sh = $grp.Item(58)
$sh.Left = 581.9903937007874
$sh.Top = 207.58582677165353
$sh = $grp.Item(59)
$sh.Left = 582.2571653543307
$sh.Top = 202.5843307086614
$sh = $grp.Item(60)
$sh.Left = 585.1913385826772
$sh.Top = 205.71858267716536
$sh = $grp.Item(61)
$sh.Left = 590.7596062992126
$sh.Top = 202.91779527559055
$sh = $grp.Item(62)
$sh.Left = 595.727716535433
$sh.Top = 203.95141732283466
$sh = $grp.Item(63)
$sh.Left = 591.4931496062992
$sh.Top = 204.05141732283465
$sh = $grp.Item(64)
$sh.Left = 599.1954330708661
$sh.Top = 202.5843307086614
$sh = $grp.Item(65)
$sh.Left = 599.1287401574804
$sh.Top = 204.81834645669292
$sh = $grp.Item(66)
$sh.Left = 604.7970078740158
$sh.Top = 203.05110236220472
$sh = $grp.Item(67)
$sh.Left = 604.6303149606299
$sh.Top = 205.2184251968504
$sh = $grp.Item(68)
$sh.Left = 601.1959842519685
$sh.Top = 202.5843307086614
$sh = $grp.Item(69)
$sh.Left = 601.9295275590551
$sh.Top = 205.0183464566929
$sh = $grp.Item(70)
$sh.Left = 603.1632283464567
$sh.Top = 205.0183464566929
$sh = $grp.Item(71)
$sh.Left = 601.9295275590551
$sh.Top = 206.1520472440945
$sh = $grp.Item(72)
$sh.Left = 603.1632283464567
$sh.Top = 206.1520472440945
$sh = $grp.Item(73)
$sh.Left = 600.9625984251968
$sh.Top = 207.48574803149606
$sh = $grp.Item(74)
$sh.Left = 443.75023622047246
$sh.Top = 254.33385826771652
$sh = $grp.Item(75)
$sh.Left = 444.9172440944882
$sh.Top = 254.93409448818898
$sh = $grp.Item(76)
$sh.Left = 444.883937007874
$sh.Top = 256.8346456692913
$sh = $grp.Item(77)
$sh.Left = 446.6844094488189
$sh.Top = 254.40055118110237
$sh = $grp.Item(78)
$sh.Left = 452.48614173228344
$sh.Top = 254.33385826771652
$sh = $grp.Item(79)
$sh.Left = 452.3527559055118
$sh.Top = 257.10133858267716
$sh = $grp.Item(80)
$sh.Left = 454.5867716535433
$sh.Top = 254.13385826771653
$sh = $grp.Item(81)
$sh.Left = 455.6870866141732
$sh.Top = 256.73456692913385
$sh = $grp.Item(82)
$sh.Left = 455.6870866141732
$sh.Top = 257.90157480314963
$sh = $grp.Item(83)
$sh.Left = 455.6870866141732
$sh.Top = 259.06858267716535
$sh = $grp.Item(84)
$sh.Left = 460.82188976377955
$sh.Top = 254.33385826771652
$sh = $grp.Item(85)
$sh.Left = 462.28897637795274
$sh.Top = 254.90070866141733
$sh = $grp.Item(86)
$sh.Left = 462.42236220472444
$sh.Top = 256.467874015748
$sh = $grp.Item(87)
$sh.Left = 464.7563779527559
$sh.Top = 258.06826771653544
$sh = $grp.Item(88)
$sh.Left = 469.6911811023622
$sh.Top = 254.43393700787402
$sh = $grp.Item(89)
$sh.Left = 472.1251968503937
$sh.Top = 254.0671653543307
$sh = $grp.Item(90)
$sh.Left = 473.65897637795274
$sh.Top = 255.33417322834646
$sh = $grp.Item(91)
$sh.Left = 472.09188976377953
$sh.Top = 257.168031496063
$sh = $grp.Item(92)
$sh.Left = 473.29220472440943
$sh.Top = 257.66818897637796
$sh = $grp.Item(93)
$sh.Left = 473.29220472440943
$sh.Top = 258.7351968503937
$sh = $grp.Item(94)
$sh.Left = 427.36393700787403
$sh.Top = 305.88338582677164
$sh = $grp.Item(95)
$sh.Left = 428.53094488188975
$sh.Top = 306.4835433070866
$sh = $grp.Item(96)
$sh.Left = 428.49755905511813
$sh.Top = 308.38409448818896
$sh = $grp.Item(97)
$sh.Left = 430.29811023622045
$sh.Top = 305.9500787401575
$sh = $grp.Item(98)
$sh.Left = 436.099842519685
$sh.Top = 305.88338582677164
$sh = $grp.Item(99)
$sh.Left = 435.96645669291337
$sh.Top = 308.6508661417323
$sh = $grp.Item(100)
$sh.Left = 438.2003937007874
$sh.Top = 305.68330708661415
$sh = $grp.Item(101)
$sh.Left = 439.30070866141733
$sh.Top = 308.284094488189
$sh = $grp.Item(102)
$sh.Left = 439.30070866141733
$sh.Top = 309.4511023622047
$sh = $grp.Item(103)
$sh.Left = 439.30070866141733
$sh.Top = 310.61811023622045
$sh = $grp.Item(104)
$sh.Left = 444.60228346456694
$sh.Top = 305.68330708661415
$sh = $grp.Item(105)
$sh.Left = 445.1024409448819
$sh.Top = 307.65055118110234
$sh = $grp.Item(106)
$sh.Left = 445.76929133858266
$sh.Top = 308.2507086614173
$sh = $grp.Item(107)
$sh.Left = 446.1027559055118
$sh.Top = 308.5174803149606
$sh = $grp.Item(108)
$sh.Left = 449.0369291338583
$sh.Top = 310.6514173228346
$sh = $grp.Item(109)
$sh.Left = 452.938031496063
$sh.Top = 306.0167716535433
$sh = $grp.Item(110)
$sh.Left = 454.9052755905512
$sh.Top = 306.61692913385826
$sh = $grp.Item(111)
$sh.Left = 454.73858267716537
$sh.Top = 308.0506299212598
$sh = $grp.Item(112)
$sh.Left = 454.471811023622
$sh.Top = 310.8514960629921
$sh = $grp.Item(113)
$sh.Left = 457.005905511811
$sh.Top = 305.68330708661415
$sh = $grp.Item(114)
$sh.Left = 457.33937007874016
$sh.Top = 305.7166141732283
$sh = $grp.Item(115)
$sh.Left = 296.7667716535433
$sh.Top = 358.4740157480315
$sh = $grp.Item(116)
$sh.Left = 296.9668503937008
$sh.Top = 362.5085039370079
$sh = $grp.Item(117)
$sh.Left = 301.10141732283466
$sh.Top = 358.373937007874
$sh = $grp.Item(118)
$sh.Left = 301.73488188976376
$sh.Top = 362.70858267716534
$sh = $grp.Item(119)
$sh.Left = 305.30259842519683
$sh.Top = 358.6073228346457
$sh = $grp.Item(120)
$sh.Left = 306.33622047244097
$sh.Top = 359.2075590551181
$sh = $grp.Item(121)
$sh.Left = 306.33622047244097
$sh.Top = 361.1414173228346
$sh = $grp.Item(122)
$sh.Left = 307.93669291338585
$sh.Top = 358.4073228346457
$sh = $grp.Item(123)
$sh.Left = 310.0040157480315
$sh.Top = 359.64102362204727
$sh = $grp.Item(124)
$sh.Left = 309.2370866141732
$sh.Top = 360.84133858267717
$sh = $grp.Item(125)
$sh.Left = 309.2370866141732
$sh.Top = 361.87496062992125
$sh = $grp.Item(126)
$sh.Left = 313.87181102362206
$sh.Top = 358.373937007874
$sh = $grp.Item(127)
$sh.Left = 315.67228346456693
$sh.Top = 361.27480314960627
$sh = $grp.Item(128)
$sh.Left = 317.97299212598426
$sh.Top = 361.27480314960627
$sh = $grp.Item(129)
$sh.Left = 322.2408661417323
$sh.Top = 358.34062992125985
$sh = $grp.Item(130)
$sh.Left = 322.3742519685039
$sh.Top = 359.84102362204726
$sh = $grp.Item(131)
$sh.Left = 324.7416535433071
$sh.Top = 360.2411811023622
$sh = $grp.Item(132)
$sh.Left = 324.3748818897638
$sh.Top = 362.7752755905512
$sh = $grp.Item(133)
$sh.Left = 184.44622047244096
$sh.Top = 409.9567716535433
$sh = $grp.Item(134)
$sh.Left = 182.57897637795276
$sh.Top = 410.2902362204724
$sh = $grp.Item(135)
$sh.Left = 184.07944881889765
$sh.Top = 410.75700787401576
$sh = $grp.Item(136)
$sh.Left = 190.68141732283465
$sh.Top = 410.090157480315
$sh = $grp.Item(137)
$sh.Left = 199.05047244094487
$sh.Top = 409.9234645669291
$sh = $grp.Item(138)
$sh.Left = 201.41787401574803
$sh.Top = 413.1244094488189
$sh = $grp.Item(139)
$sh.Left = 201.41787401574803
$sh.Top = 414.52480314960627
$sh = $grp.Item(140)
$sh.Left = 208.01976377952755
$sh.Top = 409.99015748031496
$sh = $grp.Item(141)
$sh.Left = 207.61968503937007
$sh.Top = 411.99070866141733
$sh = $grp.Item(142)
$sh.Left = 207.5196062992126
$sh.Top = 414.158031496063
$sh = $grp.Item(143)
$sh.Left = 209.92031496062992
$sh.Top = 410.22354330708663
$sh = $grp.Item(144)
$sh.Left = 209.5535433070866
$sh.Top = 413.55787401574804
$sh = $grp.Item(145)
$sh.Left = 211.1540157480315
$sh.Top = 414.158031496063
$sh = $grp.Item(146)
$sh.Left = 216.12212598425197
$sh.Top = 409.9234645669291
$sh = $grp.Item(147)
$sh.Left = 218.48952755905512
$sh.Top = 413.1244094488189
$sh = $grp.Item(148)
$sh.Left = 218.48952755905512
$sh.Top = 414.52480314960627
$sh = $grp.Item(149)
$sh.Left = 226.79188976377952
$sh.Top = 410.02346456692914
$sh = $grp.Item(150)
$sh.Left = 224.8247244094488
$sh.Top = 410.6903149606299
$sh = $grp.Item(151)
$sh.Left = 225.45818897637795
$sh.Top = 411.32385826771656
$sh = $grp.Item(152)
$sh.Left = 233.2271653543307
$sh.Top = 410.1568503937008
$sh = $grp.Item(153)
$sh.Left = 233.09377952755906
$sh.Top = 412.92433070866144
$sh = $grp.Item(154)
$sh.Left = 235.72787401574803
$sh.Top = 409.99015748031496
$sh = $grp.Item(155)
$sh.Left = 235.89456692913384
$sh.Top = 412.5908661417323
$sh = $grp.Item(156)
$sh.Left = 158.89110236220472
$sh.Top = 125.24543307086614
$sh.Width = 544.8168503937007
$sh.Height = 319.6067716535433
$sh = $grp.Item(157)
$sh.Left = 151.57952755905512
$sh.Top = 410.8051181102362
$sh = $grp.Item(158)
$sh.Left = 152.84708661417324
$sh.Top = 413.9418897637795
$sh = $grp.Item(159)
$sh.Left = 151.52795275590552
$sh.Top = 359.345905511811
$sh = $grp.Item(160)
$sh.Left = 151.39047244094488
$sh.Top = 307.7963779527559
$sh = $grp.Item(161)
$sh.Left = 152.3744094488189
$sh.Top = 308.89771653543306
$sh = $grp.Item(162)
$sh.Left = 151.4592125984252
$sh.Top = 256.1566929133858
$sh = $grp.Item(163)
$sh.Left = 151.56228346456692
$sh.Top = 204.60724409448818
$sh = $grp.Item(164)
$sh.Left = 151.8115748031496
$sh.Top = 153.14795275590552
$sh = $grp.Item(165)
$sh.Left = 156.15141732283465
$sh.Top = 413.9225196850394
$sh = $grp.Item(166)
$sh.Left = 156.15141732283465
$sh.Top = 362.37307086614175
$sh = $grp.Item(167)
$sh.Left = 156.15141732283465
$sh.Top = 310.8235433070866
$sh = $grp.Item(168)
$sh.Left = 156.15141732283465
$sh.Top = 259.27409448818895
$sh = $grp.Item(169)
$sh.Left = 156.15141732283465
$sh.Top = 207.72456692913386
$sh = $grp.Item(170)
$sh.Left = 156.15141732283465
$sh.Top = 156.1751181102362
$sh = $grp.Item(171)
$sh.Left = 158.89110236220472
$sh.Top = 444.85220472440943
$sh = $grp.Item(172)
$sh.Left = 309.2244094488189
$sh.Top = 444.85220472440943
$sh = $grp.Item(173)
$sh.Left = 459.5577165354331
$sh.Top = 444.85220472440943
$sh = $grp.Item(174)
$sh.Left = 609.8910236220472
$sh.Top = 444.85220472440943
$sh = $grp.Item(175)
$sh.Left = 156.78779527559055
$sh.Top = 449.69346456692915
$sh = $grp.Item(176)
$sh.Left = 157.56984251968504
$sh.Top = 450.32086614173227
$sh = $grp.Item(177)
$sh.Left = 300.1064566929134
$sh.Top = 449.7837007874016
$sh = $grp.Item(178)
$sh.Left = 304.6740157480315
$sh.Top = 449.69346456692915
$sh = $grp.Item(179)
$sh.Left = 305.45606299212596
$sh.Top = 450.32086614173227
$sh = $grp.Item(180)
$sh.Left = 309.56818897637794
$sh.Top = 449.69346456692915
$sh = $grp.Item(181)
$sh.Left = 310.3502362204724
$sh.Top = 450.32086614173227
$sh = $grp.Item(182)
$sh.Left = 314.46228346456695
$sh.Top = 449.69346456692915
$sh = $grp.Item(183)
$sh.Left = 315.24433070866144
$sh.Top = 450.32086614173227
$sh = $grp.Item(184)
$sh.Left = 450.2120472440945
$sh.Top = 449.69346456692915
$sh = $grp.Item(185)
$sh.Left = 455.00732283464566
$sh.Top = 449.69346456692915
$sh = $grp.Item(186)
$sh.Left = 455.78937007874015
$sh.Top = 450.32086614173227
$sh = $grp.Item(187)
$sh.Left = 459.9014960629921
$sh.Top = 449.69346456692915
$sh = $grp.Item(188)
$sh.Left = 460.68346456692916
$sh.Top = 450.32086614173227
$sh = $grp.Item(189)
$sh.Left = 464.7955905511811
$sh.Top = 449.69346456692915
$sh = $grp.Item(190)
$sh.Left = 465.57763779527556
$sh.Top = 450.32086614173227
$sh = $grp.Item(191)
$sh.Left = 600.437874015748
$sh.Top = 449.69346456692915
$sh = $grp.Item(192)
$sh.Left = 605.3406299212598
$sh.Top = 449.69346456692915
$sh = $grp.Item(193)
$sh.Left = 606.1226771653543
$sh.Top = 450.32086614173227
$sh = $grp.Item(194)
$sh.Left = 610.2348031496063
$sh.Top = 449.69346456692915
$sh = $grp.Item(195)
$sh.Left = 611.0167716535433
$sh.Top = 450.32086614173227
$sh = $grp.Item(196)
$sh.Left = 615.1288976377953
$sh.Top = 449.69346456692915
$sh = $grp.Item(197)
$sh.Left = 615.9109448818898
$sh.Top = 450.32086614173227
